# "Generate Report for Handback" - refresh the handback-status report with
# the latest handoff/handback timestamps and mark the zh-cn/de-de rows as
# machine translated ("mt") instead of human translated ("ht").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G2").Value = "2016-08-26 16:17:20"
$wsOverview.Range("G4").Value = "2016-08-26 16:17:20"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-26 16:17:15"
$wsZhCn.Range("H4").Value = "2016-08-26 16:17:15"
$wsZhCn.Range("K2").Value = "2016-08-26 16:17:38"
$wsZhCn.Range("K4").Value = "2016-08-26 16:17:38"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-26 16:17:20"
$wsDeDe.Range("H4").Value = "2016-08-26 16:17:20"
$wsDeDe.Range("K2").Value = "2016-08-26 16:17:46"
$wsDeDe.Range("K4").Value = "2016-08-26 16:17:46"
